$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 428, shifting the Andhra-district block down by one
$ws.Rows("428").Insert()

# Update Hindi/Telugu translations for the shifted rows (keys + English text already carried by the shift)
$ws.Cells.Item(429, 3).Value = "आंध्र प्रदेश में जिला चुनें"
$ws.Cells.Item(429, 4).Value = "ఆంధ్రప్రదేశ్‌లో జిల్లాను ఎంచుకోండి"
$ws.Cells.Item(430, 3).Value = "आंध्र प्रदेश में जिलों को खोजें..."
$ws.Cells.Item(430, 4).Value = "ఆంధ్రప్రదేశ్‌లో జిల్లాలను శోధించండి..."
$ws.Cells.Item(431, 3).Value = "{query} के लिए कोई जिला नहीं मिला"
$ws.Cells.Item(431, 4).Value = "{query} కోసం జిల్లాలు కనుగొనబడలేదు"
$ws.Cells.Item(432, 3).Value = "संपत्तियां उपलब्ध हैं"
$ws.Cells.Item(432, 4).Value = "ఆస్తులు అందుబాటులో ఉన్నాయి"

# Append new district / area / resort-selection / sidebar translation rows
$ws.Cells.Item(434, 1).Value = "districts.anantapur"
$ws.Cells.Item(434, 2).Value = "Anantapur"
$ws.Cells.Item(434, 3).Value = "अनंतपुर"
$ws.Cells.Item(434, 4).Value = "అనంతపురం"

$ws.Cells.Item(435, 1).Value = "districts.chittoor"
$ws.Cells.Item(435, 2).Value = "Chittoor"
$ws.Cells.Item(435, 3).Value = "चित्तूर"
$ws.Cells.Item(435, 4).Value = "చిత్తూరు"

$ws.Cells.Item(436, 1).Value = "districts.eastgodavari"
$ws.Cells.Item(436, 2).Value = "East Godavari"
$ws.Cells.Item(436, 3).Value = "पूर्वी गोदावरी"
$ws.Cells.Item(436, 4).Value = "తూర్పు గోదావరి"

$ws.Cells.Item(437, 1).Value = "districts.guntur"
$ws.Cells.Item(437, 2).Value = "Guntur"
$ws.Cells.Item(437, 3).Value = "गुंटूर"
$ws.Cells.Item(437, 4).Value = "గుంటూరు"

$ws.Cells.Item(438, 1).Value = "districts.kadapa"
$ws.Cells.Item(438, 2).Value = "Kadapa"
$ws.Cells.Item(438, 3).Value = "कडपा"
$ws.Cells.Item(438, 4).Value = "కడప"

$ws.Cells.Item(439, 1).Value = "districts.krishna"
$ws.Cells.Item(439, 2).Value = "Krishna"
$ws.Cells.Item(439, 3).Value = "कृष्णा"
$ws.Cells.Item(439, 4).Value = "కృష్ణా"

$ws.Cells.Item(440, 1).Value = "districts.kurnool"
$ws.Cells.Item(440, 2).Value = "Kurnool"
$ws.Cells.Item(440, 3).Value = "कर्नूल"
$ws.Cells.Item(440, 4).Value = "కర్నూలు"

$ws.Cells.Item(441, 1).Value = "districts.nellore"
$ws.Cells.Item(441, 2).Value = "Nellore"
$ws.Cells.Item(441, 3).Value = "नेल्लोर"
$ws.Cells.Item(441, 4).Value = "నెల్లూరు"

$ws.Cells.Item(442, 1).Value = "districts.srikakulam"
$ws.Cells.Item(442, 2).Value = "Srikakulam"
$ws.Cells.Item(442, 3).Value = "श्रीकाकुलम"
$ws.Cells.Item(442, 4).Value = "శ్రీకాకుళం"

$ws.Cells.Item(443, 1).Value = "districts.visakhapatnam"
$ws.Cells.Item(443, 2).Value = "Visakhapatnam"
$ws.Cells.Item(443, 3).Value = "विशाखापत्तनम"
$ws.Cells.Item(443, 4).Value = "విశాఖపట్టణం"

$ws.Cells.Item(444, 1).Value = "districts.vizianagaram"
$ws.Cells.Item(444, 2).Value = "Vizianagaram"
$ws.Cells.Item(444, 3).Value = "विजयनगरम"
$ws.Cells.Item(444, 4).Value = "విజయనగరం"

$ws.Cells.Item(445, 1).Value = "districts.westgodavari"
$ws.Cells.Item(445, 2).Value = "West Godavari"
$ws.Cells.Item(445, 3).Value = "पश्चिम गोदावरी"
$ws.Cells.Item(445, 4).Value = "పశ్చిమ గోదావరి"

$ws.Cells.Item(447, 1).Value = "areas.akkayapalem"
$ws.Cells.Item(447, 2).Value = "Akkayapalem"
$ws.Cells.Item(447, 3).Value = "अक्कायापलेम"
$ws.Cells.Item(447, 4).Value = "అక్కయ్యపాలెం"

$ws.Cells.Item(448, 1).Value = "areas.anandapuram"
$ws.Cells.Item(448, 2).Value = "Anandapuram"
$ws.Cells.Item(448, 3).Value = "आनंदपुरम"
$ws.Cells.Item(448, 4).Value = "ఆనందపురం"

$ws.Cells.Item(449, 1).Value = "areas.boyapalem"
$ws.Cells.Item(449, 2).Value = "Boyapalem"
$ws.Cells.Item(449, 3).Value = "बोयापलेम"
$ws.Cells.Item(449, 4).Value = "బోయపాలెం"

$ws.Cells.Item(450, 1).Value = "areas.chinnagadili"
$ws.Cells.Item(450, 2).Value = "Chinna Gadili"
$ws.Cells.Item(450, 3).Value = "चिन्ना गडिली"
$ws.Cells.Item(450, 4).Value = "చిన్న గడిలి"

$ws.Cells.Item(451, 1).Value = "areas.dwarkanagar"
$ws.Cells.Item(451, 2).Value = "Dwarka Nagar"
$ws.Cells.Item(451, 3).Value = "द्वारका नगर"
$ws.Cells.Item(451, 4).Value = "ద్వారకా నగర్"

$ws.Cells.Item(452, 1).Value = "areas.gajuwaka"
$ws.Cells.Item(452, 2).Value = "Gajuwaka"
$ws.Cells.Item(452, 3).Value = "गाजुवाका"
$ws.Cells.Item(452, 4).Value = "గాజువాక"

$ws.Cells.Item(453, 1).Value = "areas.kommadi"
$ws.Cells.Item(453, 2).Value = "Kommadi"
$ws.Cells.Item(453, 3).Value = "कोम्मडी"
$ws.Cells.Item(453, 4).Value = "కొమ్మడి"

$ws.Cells.Item(455, 1).Value = "selectSite.title"
$ws.Cells.Item(455, 2).Value = "Select Resort"
$ws.Cells.Item(455, 3).Value = "रिसॉर्ट चुनें"
$ws.Cells.Item(455, 4).Value = "రిసార్ట్‌ను ఎంచుకోండి"

$ws.Cells.Item(456, 1).Value = "selectSite.searchPlaceholder"
$ws.Cells.Item(456, 2).Value = "Search properties in"
$ws.Cells.Item(456, 3).Value = "में संपत्तियों को खोजें"
$ws.Cells.Item(456, 4).Value = "లో ఆస్తులను శోధించండి"

$ws.Cells.Item(457, 1).Value = "selectSite.noResults"
$ws.Cells.Item(457, 2).Value = "No properties found"
$ws.Cells.Item(457, 3).Value = "कोई संपत्ति नहीं मिली"
$ws.Cells.Item(457, 4).Value = "ఆస్తులు కనుగొనబడలేదు"

$ws.Cells.Item(459, 1).Value = "sidebar_menu_my_properties"
$ws.Cells.Item(459, 2).Value = "My Properties"
$ws.Cells.Item(459, 3).Value = "मेरी प्रॉपर्टीज़"
$ws.Cells.Item(459, 4).Value = "నా ఆస్తులు"

# Update selection to mirror the final authored view state
$ws.Range("J465").Select() | Out-Null
